{"js": "// Update the worksheet date and all 100 math-fact cells in the table.\n// Each [oldText, newText] pair is applied with body.search()/insertText()\n// so the edit is anchored to the literal text rather than a fragile\n// paragraph/cell index.\nconst replacements = [\n  [\"2025-10-20 Monday\", \"2025-10-21 Tuesday\"],\n  [\"83-24=\", \"53-34=\"],\n  [\"94+1=\", \"67+1=\"],\n  [\"14-5=\", \"71-0=\"],\n  [\"15+74=\", \"54+36=\"],\n  [\"45+4=\", \"61-37=\"],\n  [\"3+13=\", \"47-28=\"],\n  [\"21+0=\", \"96-85=\"],\n  [\"62+3=\", \"59-9=\"],\n  [\"52-5=\", \"82-81=\"],\n  [\"8+37=\", \"53-0=\"],\n  [\"38+50=\", \"41+40=\"],\n  [\"80+1=\", \"45-9=\"],\n  [\"80+12=\", \"50+44=\"],\n  [\"28+26=\", \"94-4=\"],\n  [\"31+58=\", \"52-34=\"],\n  [\"49-34=\", \"39+55=\"],\n  [\"81+2=\", \"61-0=\"],\n  [\"50-44=\", \"10+72=\"],\n  [\"60+13=\", \"56+26=\"],\n  [\"92+5=\", \"69+26=\"],\n  [\"57-38=\", \"75-1=\"],\n  [\"34+7=\", \"56-29=\"],\n  [\"40-19=\", \"26+6=\"],\n  [\"99-47=\", \"84-69=\"],\n  [\"19+33=\", \"94-36=\"],\n  [\"93-39=\", \"68-17=\"],\n  [\"78-77=\", \"17-12=\"],\n  [\"60-32=\", \"10+83=\"],\n  [\"17+34=\", \"28+67=\"],\n  [\"38+32=\", \"12+75=\"],\n  [\"92-16=\", \"0+31=\"],\n  [\"90-43=\", \"18+33=\"],\n  [\"41+47=\", \"95+4=\"],\n  [\"14+43=\", \"26+0=\"],\n  [\"77+21=\", \"86-21=\"],\n  [\"7-6=\", \"29+38=\"],\n  [\"67-47=\", \"98-56=\"],\n  [\"10+82=\", \"56-25=\"],\n  [\"85+7=\", \"11-9=\"],\n  [\"88-29=\", \"76-45=\"],\n  [\"25+45=\", \"67-54=\"],\n  [\"40+7=\", \"23+1=\"],\n  [\"5+58=\", \"27+17=\"],\n  [\"42-22=\", \"96-84=\"],\n  [\"19+51=\", \"24+30=\"],\n  [\"78-51=\", \"61+35=\"],\n  [\"69-19=\", \"19-6=\"],\n  [\"34+0=\", \"0+11=\"],\n  [\"39-0=\", \"95-15=\"],\n  [\"95-48=\", \"44-16=\"],\n  [\"47-17=\", \"10+2=\"],\n  [\"88-7=\", \"62-6=\"],\n  [\"68+10=\", \"50+37=\"],\n  [\"10+64=\", \"21+8=\"],\n  [\"66-9=\", \"28-25=\"],\n  [\"65-46=\", \"37+62=\"],\n  [\"80-78=\", \"77-49=\"],\n  [\"36-27=\", \"70+5=\"],\n  [\"45+44=\", \"97-22=\"],\n  [\"6+80=\", \"64-47=\"],\n  [\"58+30=\", \"91-74=\"],\n  [\"96-39=\", \"0+22=\"],\n  [\"0+18=\", \"76+4=\"],\n  [\"28+13=\", \"70-4=\"],\n  [\"85-14=\", \"24+29=\"],\n  [\"34+24=\", \"84-78=\"],\n  [\"3+6=\", \"46-17=\"],\n  [\"26+55=\", \"96-33=\"],\n  [\"30+29=\", \"75-47=\"],\n  [\"56-51=\", \"39+8=\"],\n  [\"20+73=\", \"2+28=\"],\n  [\"22+62=\", \"20-4=\"],\n  [\"88-45=\", \"58-24=\"],\n  [\"45-45=\", \"18+12=\"],\n  [\"51+34=\", \"55+14=\"],\n  [\"3+66=\", \"34-20=\"],\n  [\"3+79=\", \"30-16=\"],\n  [\"32-15=\", \"77-13=\"],\n  [\"97-20=\", \"2+55=\"],\n  [\"24+47=\", \"80-37=\"],\n  [\"14+85=\", \"42+42=\"],\n  [\"32+12=\", \"75-7=\"],\n  [\"99-61=\", \"54-35=\"],\n  [\"32+63=\", \"45+12=\"],\n  [\"77-3=\", \"17+61=\"],\n  [\"63-21=\", \"38+40=\"],\n  [\"97-4=\", \"50-21=\"],\n  [\"16+71=\", \"12+67=\"],\n  [\"73+24=\", \"43+54=\"],\n  [\"86-48=\", \"43+21=\"],\n  [\"59+12=\", \"87-17=\"],\n  [\"60+12=\", \"71-59=\"],\n  [\"88-18=\", \"10+60=\"],\n  [\"38-15=\", \"53-0=\"],\n  [\"60-12=\", \"99-59=\"],\n  [\"76-42=\", \"75-54=\"],\n  [\"66-12=\", \"74+0=\"],\n  [\"58-25=\", \"5+36=\"],\n  [\"86-41=\", \"48-1=\"],\n  [\"31+68=\", \"82-18=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 100 math-fact cells in the table.\n# Each (old, new) pair below is applied with Find/Replace so the edit\n# is anchored to the literal text rather than a fragile cell index.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2025-10-20 Monday', '2025-10-21 Tuesday'),\n    @('83-24=', '53-34='),\n    @('94+1=', '67+1='),\n    @('14-5=', '71-0='),\n    @('15+74=', '54+36='),\n    @('45+4=', '61-37='),\n    @('3+13=', '47-28='),\n    @('21+0=', '96-85='),\n    @('62+3=', '59-9='),\n    @('52-5=', '82-81='),\n    @('8+37=', '53-0='),\n    @('38+50=', '41+40='),\n    @('80+1=', '45-9='),\n    @('80+12=', '50+44='),\n    @('28+26=', '94-4='),\n    @('31+58=', '52-34='),\n    @('49-34=', '39+55='),\n    @('81+2=', '61-0='),\n    @('50-44=', '10+72='),\n    @('60+13=', '56+26='),\n    @('92+5=', '69+26='),\n    @('57-38=', '75-1='),\n    @('34+7=', '56-29='),\n    @('40-19=', '26+6='),\n    @('99-47=', '84-69='),\n    @('19+33=', '94-36='),\n    @('93-39=', '68-17='),\n    @('78-77=', '17-12='),\n    @('60-32=', '10+83='),\n    @('17+34=', '28+67='),\n    @('38+32=', '12+75='),\n    @('92-16=', '0+31='),\n    @('90-43=', '18+33='),\n    @('41+47=', '95+4='),\n    @('14+43=', '26+0='),\n    @('77+21=', '86-21='),\n    @('7-6=', '29+38='),\n    @('67-47=', '98-56='),\n    @('10+82=', '56-25='),\n    @('85+7=', '11-9='),\n    @('88-29=', '76-45='),\n    @('25+45=', '67-54='),\n    @('40+7=', '23+1='),\n    @('5+58=', '27+17='),\n    @('42-22=', '96-84='),\n    @('19+51=', '24+30='),\n    @('78-51=', '61+35='),\n    @('69-19=', '19-6='),\n    @('34+0=', '0+11='),\n    @('39-0=', '95-15='),\n    @('95-48=', '44-16='),\n    @('47-17=', '10+2='),\n    @('88-7=', '62-6='),\n    @('68+10=', '50+37='),\n    @('10+64=', '21+8='),\n    @('66-9=', '28-25='),\n    @('65-46=', '37+62='),\n    @('80-78=', '77-49='),\n    @('36-27=', '70+5='),\n    @('45+44=', '97-22='),\n    @('6+80=', '64-47='),\n    @('58+30=', '91-74='),\n    @('96-39=', '0+22='),\n    @('0+18=', '76+4='),\n    @('28+13=', '70-4='),\n    @('85-14=', '24+29='),\n    @('34+24=', '84-78='),\n    @('3+6=', '46-17='),\n    @('26+55=', '96-33='),\n    @('30+29=', '75-47='),\n    @('56-51=', '39+8='),\n    @('20+73=', '2+28='),\n    @('22+62=', '20-4='),\n    @('88-45=', '58-24='),\n    @('45-45=', '18+12='),\n    @('51+34=', '55+14='),\n    @('3+66=', '34-20='),\n    @('3+79=', '30-16='),\n    @('32-15=', '77-13='),\n    @('97-20=', '2+55='),\n    @('24+47=', '80-37='),\n    @('14+85=', '42+42='),\n    @('32+12=', '75-7='),\n    @('99-61=', '54-35='),\n    @('32+63=', '45+12='),\n    @('77-3=', '17+61='),\n    @('63-21=', '38+40='),\n    @('97-4=', '50-21='),\n    @('16+71=', '12+67='),\n    @('73+24=', '43+54='),\n    @('86-48=', '43+21='),\n    @('59+12=', '87-17='),\n    @('60+12=', '71-59='),\n    @('88-18=', '10+60='),\n    @('38-15=', '53-0='),\n    @('60-12=', '99-59='),\n    @('76-42=', '75-54='),\n    @('66-12=', '74+0='),\n    @('58-25=', '5+36='),\n    @('86-41=', '48-1='),\n    @('31+68=', '82-18='),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    [void]$find.Execute(\n        $pair[0],  # FindText\n        $false,    # MatchCase\n        $true,     # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $pair[1],  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n}\n"}
